$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2026315789473684
$ws.Range("C2").Value = 0.5486842105263158
$ws.Range("J2").Value = 0.01578947368421053
$ws.Range("P2").Value = 0.1526315789473684
$ws.Range("S2").Value = 0.08026315789473684

$ws.Range("B3").Value = 0.009324009324009324
$ws.Range("C3").Value = 0.02097902097902098
$ws.Range("J3").Value = 0.0372960372960373
$ws.Range("P3").Value = 0.7435897435897436
$ws.Range("S3").Value = 0.1888111888111888

$ws.Range("J4").Value = 0.02586206896551724
$ws.Range("P4").Value = 0.6637931034482759
$ws.Range("S4").Value = 0.3103448275862069

$ws.Range("P5").Value = 0.6666666666666666
$ws.Range("S5").Value = 0.3333333333333333

$ws.Range("B6").Value = 0.08171206225680934
$ws.Range("D6").Value = 0.003891050583657588
$ws.Range("F6").Value = 0.06031128404669261
$ws.Range("J6").Value = 0.2684824902723735
$ws.Range("O6").Value = 0.01167315175097276
$ws.Range("Q6").Value = 0.1653696498054475
$ws.Range("R6").Value = 0.08171206225680934
$ws.Range("S6").Value = 0.3268482490272374

$ws.Range("B7").Value = 0.1090487238979118
$ws.Range("D7").Value = 0.0185614849187935
$ws.Range("F7").Value = 0.04872389791183294
$ws.Range("J7").Value = 0.1345707656612529
$ws.Range("O7").Value = 0.02320185614849188
$ws.Range("Q7").Value = 0.1763341067285383
$ws.Range("R7").Value = 0.08816705336426914
$ws.Range("S7").Value = 0.4013921113689095

$ws.Range("B8").Value = 0.09749552772808587
$ws.Range("D8").Value = 0.01520572450805009
$ws.Range("E8").Value = 0.0008944543828264759
$ws.Range("F8").Value = 0.05724508050089445
$ws.Range("J8").Value = 0.1198568872987478
$ws.Range("O8").Value = 0.01788908765652952
$ws.Range("Q8").Value = 0.1887298747763864
$ws.Range("R8").Value = 0.09570661896243292
$ws.Range("S8").Value = 0.4069767441860465

$ws.Range("B9").Value = 0.1030042918454936
$ws.Range("D9").Value = 0.02145922746781116
$ws.Range("E9").Value = 0.004291845493562232
$ws.Range("F9").Value = 0.04721030042918455
$ws.Range("J9").Value = 0.1094420600858369
$ws.Range("O9").Value = 0.02575107296137339
$ws.Range("Q9").Value = 0.1995708154506438
$ws.Range("R9").Value = 0.09656652360515021
$ws.Range("S9").Value = 0.3927038626609442

$ws.Range("B10").Value = 0.1128304319793682
$ws.Range("D10").Value = 0.02482269503546099
$ws.Range("E10").Value = 0.001289490651192779
$ws.Range("F10").Value = 0.06769825918762089
$ws.Range("J10").Value = 0.1054158607350097
$ws.Range("O10").Value = 0.01579626047711154
$ws.Range("Q10").Value = 0.2176015473887814
$ws.Range("R10").Value = 0.08381689232753063
$ws.Range("S10").Value = 0.3707285622179239

$ws.Range("G11").Value = 0.1307086614173228
$ws.Range("J11").Value = 0.08818897637795275
$ws.Range("K11").Value = 0.1716535433070866
$ws.Range("L11").Value = 0.5984251968503937
$ws.Range("S11").Value = 0.01102362204724409

$ws.Range("G12").Value = 0.6990049751243781
$ws.Range("J12").Value = 0.2238805970149254
$ws.Range("K12").Value = 0.009950248756218905
$ws.Range("L12").Value = 0.0472636815920398
$ws.Range("S12").Value = 0.01990049751243781

$ws.Range("G13").Value = 0.7090909090909091
$ws.Range("J13").Value = 0.2272727272727273
$ws.Range("S13").Value = 0.06363636363636363

$ws.Range("F15").Value = 0.01834862385321101
$ws.Range("H15").Value = 0.1614678899082569
$ws.Range("I15").Value = 0.06238532110091743
$ws.Range("J15").Value = 0.344954128440367
$ws.Range("K15").Value = 0.06238532110091743
$ws.Range("M15").Value = 0.02018348623853211
$ws.Range("N15").Value = 0.003669724770642202
$ws.Range("O15").Value = 0.08256880733944955
$ws.Range("S15").Value = 0.2440366972477064

$ws.Range("F16").Value = 0.01202404809619238
$ws.Range("H16").Value = 0.1703406813627255
$ws.Range("I16").Value = 0.08016032064128256
$ws.Range("J16").Value = 0.4128256513026052
$ws.Range("K16").Value = 0.1142284569138277
$ws.Range("M16").Value = 0.01402805611222445
$ws.Range("O16").Value = 0.04609218436873747
$ws.Range("S16").Value = 0.1503006012024048

$ws.Range("F17").Value = 0.02298850574712644
$ws.Range("H17").Value = 0.1847922192749779
$ws.Range("I17").Value = 0.08930150309460655
$ws.Range("J17").Value = 0.4438549955791335
$ws.Range("K17").Value = 0.08399646330680813
$ws.Range("M17").Value = 0.01326259946949602
$ws.Range("N17").Value = 0.0008841732979664014
$ws.Range("O17").Value = 0.04951370468611848
$ws.Range("S17").Value = 0.1114058355437666

$ws.Range("F18").Value = 0.02443991853360489
$ws.Range("H18").Value = 0.1995926680244399
$ws.Range("I18").Value = 0.09164969450101833
$ws.Range("J18").Value = 0.3971486761710794
$ws.Range("K18").Value = 0.07942973523421588
$ws.Range("M18").Value = 0.01425661914460285
$ws.Range("O18").Value = 0.05906313645621181
$ws.Range("S18").Value = 0.1344195519348269

$ws.Range("F19").Value = 0.01561991539212496
$ws.Range("H19").Value = 0.2111942726976896
$ws.Range("I19").Value = 0.08005206638464042
$ws.Range("J19").Value = 0.3686950862349496
$ws.Range("K19").Value = 0.09729905629677839
$ws.Range("M19").Value = 0.02440611780019525
$ws.Range("N19").Value = 0.0009762447120078099
$ws.Range("O19").Value = 0.07679791734461439
$ws.Range("S19").Value = 0.1249593231369997
